# Update the "target_data" sheet's Table1[target_range_end] column (column G)
# Rows 2-13: single-row ranges that extended to column AX now extend to the
#            very last column, XFD (B1:AX1 -> B1:XFD1, etc.)
# Rows 14-25: single-column ranges that extended to row 50 now extend to
#            row 163840 (A2:A50 -> A2:A163840, etc.)

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("target_data")

for ($i = 0; $i -lt 12; $i++) {
    $row = 2 + $i
    $ws.Cells.Item($row, 7).Value = "XFD" + ($row - 1)
}

$colLetters = @("A", "A", "B", "C", "D", "E", "F", "G", "H", "I", "J", "K")
for ($i = 0; $i -lt $colLetters.Length; $i++) {
    $row = 14 + $i
    $ws.Cells.Item($row, 7).Value = $colLetters[$i] + "163840"
}
